# Fruta / hortaliza, semanal
# Insert a new weekly observation at row 416 of "Sheet1", pushing the
# existing rows 416-430 down to 417-431 (the sheet's dimension grows
# from A1:R430 to A1:R431), then populate the newly inserted row with
# the new week's data (2023-08-09 / Excel serial 45147).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 416, shifting rows 416:430 down to 417:431.
$ws.Rows("416:416").Insert()

# Fill in the new row 416 with this week's record.
$ws.Cells.Item(416, 1).Value = 8
$ws.Cells.Item(416, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(416, 3).Value = "Coquimbo"
$ws.Cells.Item(416, 4).Value = 45147
$ws.Cells.Item(416, 5).Value = 4
$ws.Cells.Item(416, 6).Value = 100112031
$ws.Cells.Item(416, 7).Value = "Poroto verde"
$ws.Cells.Item(416, 8).Value = "Magnum"
$ws.Cells.Item(416, 9).Value = "Primera"
$ws.Cells.Item(416, 10).Value = 400
$ws.Cells.Item(416, 11).Value = 32000
$ws.Cells.Item(416, 12).Value = 33000
$ws.Cells.Item(416, 13).Value = 32500
$ws.Cells.Item(416, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(416, 15).Value = "Perú"
$ws.Cells.Item(416, 16).Value = 1300
$ws.Cells.Item(416, 17).Value = 25
$ws.Cells.Item(416, 18).Value = "Hortaliza"
